# Deploy til rot: LMDI - 2beda24b52b3078571421cf68055f15e4eb1cf64
#
# 1. Remove the "Include #1" worksheet (sheetId 3 / sheet3.xml) entirely.
# 2. Bump the Version value in the Metadata sheet from 0.9.4 to 1.0.0.
# 3. Bump the Date value in the Metadata sheet to the new publish date.

$wb = $excel.ActiveWorkbook

# 1. Delete the "Include #1" worksheet.
$wsInclude1 = $wb.Worksheets.Item("Include #1")
$wsInclude1.Delete()

# 2 & 3. Update the Version and Date values on the Metadata sheet.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.0.0"
$wsMeta.Range("B8").Value = "2025-03-18T14:32:32+00:00"
